$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns for a.md / b.md rows (en-US "status" mirror) ---
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
# Status column for a.md / b.md rows
$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("B3").Value = $newStatus

# Populate "Latest Target File" (E) and "Latest Handback File" (F) now that handback
# has happened - mirror the source (A) and handoff (C) file links/values.
$zhcn.Range("E2").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/6b2327e92806b8cda21b8756b83d118209cd5163/e2e/a.md", "", "", "a.md") | Out-Null

$zhcn.Range("F2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78d4b0604ea07b7618673fbe7abd0ad0d9f8688e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

$zhcn.Range("E3").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/6b2327e92806b8cda21b8756b83d118209cd5163/e2e/a.md", "", "", "a.md") | Out-Null

$zhcn.Range("F3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78d4b0604ea07b7618673fbe7abd0ad0d9f8688e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

# Latest Handback DateTime (G) now stamped with the handback time
$zhcn.Range("G2").Value = "2016-02-29 03:54:11"
$zhcn.Range("G3").Value = "2016-02-29 03:54:11"

# --- de-de sheet ---
$dede.Range("B2").Value = $newStatus
$dede.Range("B3").Value = $newStatus

$dede.Range("E2").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/6b2327e92806b8cda21b8756b83d118209cd5163/e2e/a.md", "", "", "a.md") | Out-Null

$dede.Range("F2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e7ab39da84844b848f0aef1f29b3b1b84daa955/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

$dede.Range("E3").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/6b2327e92806b8cda21b8756b83d118209cd5163/e2e/a.md", "", "", "a.md") | Out-Null

$dede.Range("F3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e7ab39da84844b848f0aef1f29b3b1b84daa955/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

$dede.Range("G2").Value = "2016-02-29 03:54:34"
$dede.Range("G3").Value = "2016-02-29 03:54:34"
